# Fixing the huge difference between the hourly and quarterly forecasted values
#
# Column A (timestamps, rows 2-97) are shifted forward by 3 days (the serial
# date values increase by exactly 3.0), and column B (Notified Production MW,
# rows 22-86) is updated with corrected forecast values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift every timestamp in column A (rows 2 through 97) by +3 days ---
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 3
}

# --- Update the quarterly forecast values in column B (rows 22 through 86) ---
$newB = @(12,9,2,2,54,49,54,62,267,281,301,321,573,598,626,651,858,881,905,927,1062,1082,1098,1113,1176,1188,1199,1207,1198,1198,1197,1194,1124,1117,1109,1096,1001,983,966,949,779,755,731,708,480,454,431,411,203,187,169,155,21,16,14,13,1,1,1,1,1,1,1,1,1)

$startRow = 22
for ($i = 0; $i -lt $newB.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newB[$i]
}
